# Auto-generated edit script applying the crypto price/volume update diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "69.884.03"
Set-TextCell "E2" "  -0.15%  "

Set-TextCell "D3" "3.931.03"
Set-TextCell "E3" "  +5.83%  "

Set-TextCell "D4" "1.00"
Set-TextCell "E4" "  -0.02%  "

Set-TextCell "D5" "611.39"
Set-TextCell "E5" "  -1.16%  "

Set-TextCell "D6" "175.90"
Set-TextCell "E6" "  -1.85%  "

Set-TextCell "D7" "3.933.44"
Set-TextCell "E7" "  +5.90%  "

Set-TextCell "E8" "  +0.00%  "

Set-TextCell "E9" "  -0.83%  "

Set-TextCell "E10" "  +1.42%  "

Set-TextCell "D11" "6.45"
Set-TextCell "E11" "  +2.34%  "

Set-TextCell "D12" "0.485"
Set-TextCell "E12" "  +0.52%  "

Set-TextCell "E13" "  -0.51%  "

Set-TextCell "E14" "  -0.11%  "

Set-TextCell "D15" "4.571.62"
Set-TextCell "E15" "  +5.44%  "

Set-TextCell "D16" "3.923.07"
Set-TextCell "E16" "  +5.65%  "

Set-TextCell "D17" "69.140.49"
Set-TextCell "E17" "  -1.24%  "

Set-TextCell "D18" "7.53"
Set-TextCell "E18" "  -0.79%  "

Set-TextCell "E19" "  -3.04%  "

Set-TextCell "D20" "16.74"
Set-TextCell "E20" "  +1.15%  "

Set-TextCell "D21" "508.65"
Set-TextCell "E21" "  +0.94%  "

Set-TextCell "D22" "9.72"
Set-TextCell "E22" "  +5.50%  "

Set-TextCell "D23" "0.753"
Set-TextCell "E23" "  +4.90%  "

Set-TextCell "B24" "Litecoin"
Set-TextCell "C24" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D24" "87.01"
Set-TextCell "E24" "  +0.89%  "

Set-TextCell "B25" "Fetch.AI"
Set-TextCell "C25" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D25" "2.46"
Set-TextCell "E25" "  -3.72%  "

Set-TextCell "D26" "0.0000141"
Set-TextCell "E26" "  +5.91%  "

Set-TextCell "D27" "12.71"
Set-TextCell "E27" "  -2.42%  "

Set-TextCell "D28" "10.48"
Set-TextCell "E28" "  -8.02%  "

Set-TextCell "E29" "  +0.31%  "

Set-TextCell "D30" "2.59"
Set-TextCell "E30" "  +4.67%  "

Set-TextCell "E31" "  +2.66%  "

Set-TextCell "D32" "33.71"
Set-TextCell "E32" "  +9.67%  "

Set-TextCell "D33" "7.96"
Set-TextCell "E33" "  +1.11%  "

Set-TextCell "E34" "  -0.27%  "

Set-TextCell "D35" "1.00"
Set-TextCell "E35" "  -0.08%  "

Set-TextCell "D36" "1.06"
Set-TextCell "E36" "  +0.93%  "

Set-TextCell "D37" "6.18"
Set-TextCell "E37" "  +1.05%  "

Set-TextCell "D38" "0.141"
Set-TextCell "E38" "  +2.26%  "

Set-TextCell "D39" "471.80"
Set-TextCell "E39" "  +9.41%  "

Set-TextCell "D40" "0.337"
Set-TextCell "E40" "  -0.91%  "

Set-TextCell "D41" "2.06"
Set-TextCell "E41" "  -0.63%  "

Set-TextCell "D42" "49.89"
Set-TextCell "E42" "  +0.00%  "

Set-TextCell "D43" "2.94"
Set-TextCell "E43" "  +1.47%  "

Set-TextCell "B44" "Cosmos"
Set-TextCell "C44" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D44" "8.61"
Set-TextCell "E44" "  -0.18%  "

Set-TextCell "B45" "Arweave"
Set-TextCell "C45" "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextCell "D45" "42.84"
Set-TextCell "E45" "  -6.25%  "

Set-TextCell "D46" "2.964.81"
Set-TextCell "E46" "  +0.02%  "

Set-TextCell "D47" "0.0366"
Set-TextCell "E47" "  +1.27%  "

Set-TextCell "B48" "InjectiveProtocol"
Set-TextCell "C48" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D48" "27.52"
Set-TextCell "E48" "  +0.35%  "

Set-TextCell "B49" "Monero"
Set-TextCell "C49" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D49" "140.05"
Set-TextCell "E49" "  +2.64%  "

Set-TextCell "E50" "  +0.01%  "

Set-TextCell "D51" "2.43"
Set-TextCell "E51" "  -2.50%  "
